$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-06-11 Wednesday" "2025-06-12 Thursday"

Replace-Text "182÷5=" "363÷9="
Replace-Text "434÷6=" "568÷3="
Replace-Text "193÷9=" "763÷7="
Replace-Text "338÷5=" "150÷4="
Replace-Text "741÷8=" "953÷3="
Replace-Text "825÷6=" "932÷8="
Replace-Text "228÷6=" "845÷4="
Replace-Text "523÷8=" "399÷5="
Replace-Text "898÷4=" "378÷6="
Replace-Text "390÷3=" "692÷8="
Replace-Text "526÷2=" "972÷4="
Replace-Text "230÷6=" "119÷8="
Replace-Text "218÷3=" "847÷7="
Replace-Text "609÷8=" "838÷8="
Replace-Text "830÷6=" "937÷7="
Replace-Text "839÷4=" "687÷2="
Replace-Text "549÷8=" "926÷8="
Replace-Text "171÷3=" "160÷3="
Replace-Text "624÷6=" "470÷3="
Replace-Text "593÷6=" "212÷5="
Replace-Text "764÷2=" "585÷2="
Replace-Text "816÷2=" "805÷3="
Replace-Text "506÷7=" "668÷3="
Replace-Text "827÷5=" "586÷4="
Replace-Text "102÷7=" "258÷4="

Write-Output "Done"
